$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 80.125
$ws.Range("I42").Value = 52
$ws.Range("J42").Value = 108.25
$ws.Range("K42").Value = 156
$ws.Range("L42").Value = 324.75
$ws.Range("M42").Value = 74
$ws.Range("N42").Value = -784.75
$ws.Range("H64").Value = 3413.9
$ws.Range("I64").Value = 3391.2856
$ws.Range("J64").Value = 3466.6667
$ws.Range("K64").Value = 3391.2856
$ws.Range("L64").Value = 3466.6667
$ws.Range("M64").Value = -3143.2856
$ws.Range("N64").Value = -3962.6667
$ws.Range("H67").Value = 3413.9
$ws.Range("I67").Value = 3391.2856
$ws.Range("J67").Value = 3466.6667
$ws.Range("K67").Value = 3391.2856
$ws.Range("L67").Value = 3466.6667
$ws.Range("M67").Value = -2533.2856
$ws.Range("N67").Value = -5182.6667
$ws.Range("H76").Value = 8786.6
$ws.Range("I76").Value = 8075
$ws.Range("J76").Value = 9045.362999999999
$ws.Range("K76").Value = 8075
$ws.Range("L76").Value = 9045.362999999999
$ws.Range("M76").Value = -7760
$ws.Range("N76").Value = -9675.362999999999
$ws.Range("H79").Value = 8786.6
$ws.Range("I79").Value = 8075
$ws.Range("J79").Value = 9045.362999999999
$ws.Range("K79").Value = 8075
$ws.Range("L79").Value = 9045.362999999999
$ws.Range("M79").Value = -6983
$ws.Range("N79").Value = -11229.363
$ws.Range("H138").Value = 2733.4
$ws.Range("I138").Value = 2342.9092
$ws.Range("J138").Value = 2959.4736
$ws.Range("K138").Value = 7028.7276
$ws.Range("L138").Value = 8878.4208
$ws.Range("M138").Value = -1888.7276
$ws.Range("N138").Value = -19158.4208

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 947.85
$ws.Range("I2").Value = 999.13336
$ws.Range("J2").Value = 794
$ws.Range("K2").Value = 999.13336
$ws.Range("L2").Value = 794
$ws.Range("M2").Value = -886.13336
$ws.Range("N2").Value = -1020
$ws.Range("H19").Value = 10000
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()
$ws.Range("H74").Value = 12197702
$ws.Range("I74").Value = 18519872
$ws.Range("J74").Value = 4945.0713
$ws.Range("K74").Value = 18519872
$ws.Range("L74").Value = 4945.0713
$ws.Range("M74").Value = -18518998
$ws.Range("N74").Value = -6693.0713
$ws.Range("H77").Value = 12197702
$ws.Range("I77").Value = 18519872
$ws.Range("J77").Value = 4945.0713
$ws.Range("K77").Value = 92599360
$ws.Range("L77").Value = 24725.3565
$ws.Range("M77").Value = -92594992
$ws.Range("N77").Value = -33461.35649999999
$ws.Range("H102").Value = 1342
$ws.Range("I102").Value = 1177.5
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1177.5
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = 444.5
$ws.Range("N102").Value = -5244
$ws.Range("H116").Value = 947.85
$ws.Range("I116").Value = 999.13336
$ws.Range("J116").Value = 794
$ws.Range("K116").Value = 999.13336
$ws.Range("L116").Value = 794
$ws.Range("M116").Value = 1294.86664
$ws.Range("N116").Value = -5382

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 947.85
$ws.Range("I3").Value = 999.13336
$ws.Range("J3").Value = 794
$ws.Range("K3").Value = 999.13336
$ws.Range("L3").Value = 794
$ws.Range("M3").Value = -885.13336
$ws.Range("N3").Value = -1022
$ws.Range("H107").Value = 2303.1428
$ws.Range("I107").Value = 2503.6667
$ws.Range("J107").Value = 1100
$ws.Range("K107").Value = 2503.6667
$ws.Range("L107").Value = 1100
$ws.Range("M107").Value = -583.6667000000002
$ws.Range("N107").Value = -4940

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 26558000
$ws.Range("I4").Value = 29421178
$ws.Range("J4").Value = 10333333
$ws.Range("K4").Value = 29421178
$ws.Range("L4").Value = 10333333
$ws.Range("M4").Value = -29421066
$ws.Range("N4").Value = -10333557
$ws.Range("H31").Value = 20841314
$ws.Range("I31").Value = 8976.571
$ws.Range("J31").Value = 166667680
$ws.Range("K31").Value = 8976.571
$ws.Range("L31").Value = 166667680
$ws.Range("M31").Value = -8681.571
$ws.Range("N31").Value = -166668270
$ws.Range("H34").Value = 20841314
$ws.Range("I34").Value = 8976.571
$ws.Range("J34").Value = 166667680
$ws.Range("K34").Value = 8976.571
$ws.Range("L34").Value = 166667680
$ws.Range("M34").Value = -8774.571
$ws.Range("N34").Value = -166668084
$ws.Range("H99").Value = 1925.3334
$ws.Range("I99").Value = 1236.25
$ws.Range("J99").Value = 3303.5
$ws.Range("K99").Value = 1236.25
$ws.Range("L99").Value = 3303.5
$ws.Range("M99").Value = 261.75
$ws.Range("N99").Value = -6299.5
$ws.Range("H107").Value = 541.7646999999999
$ws.Range("I107").Value = 655.8182
$ws.Range("J107").Value = 332.66666
$ws.Range("K107").Value = 655.8182
$ws.Range("L107").Value = 332.66666
$ws.Range("M107").Value = 1264.1818
$ws.Range("N107").Value = -4172.66666
$ws.Range("H126").Value = 1925.3334
$ws.Range("I126").Value = 1236.25
$ws.Range("J126").Value = 3303.5
$ws.Range("K126").Value = 3708.75
$ws.Range("L126").Value = 9910.5
$ws.Range("M126").Value = -1238.75
$ws.Range("N126").Value = -14850.5
$ws.Range("H140").Value = 32723.25
$ws.Range("J140").Value = 32723.25
$ws.Range("L140").Value = 32723.25
$ws.Range("N140").Value = -43083.25

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5002.6665
$ws.Range("I3").Value = 3509.375
$ws.Range("J3").Value = 6709.2856
$ws.Range("K3").Value = 10528.125
$ws.Range("L3").Value = 20127.8568
$ws.Range("M3").Value = -10416.125
$ws.Range("N3").Value = -20351.8568
$ws.Range("H82").Value = 5629
$ws.Range("I82").Value = 890
$ws.Range("J82").Value = 6358.077
$ws.Range("K82").Value = 2670
$ws.Range("L82").Value = 19074.231
$ws.Range("M82").Value = -2264
$ws.Range("N82").Value = -19886.231
$ws.Range("H85").Value = 5629
$ws.Range("I85").Value = 890
$ws.Range("J85").Value = 6358.077
$ws.Range("K85").Value = 2670
$ws.Range("L85").Value = 19074.231
$ws.Range("M85").Value = -1266
$ws.Range("N85").Value = -21882.231
$ws.Range("H118").Value = 2262.5
$ws.Range("I118").Value = 2836
$ws.Range("J118").Value = 1306.6666
$ws.Range("K118").Value = 8508
$ws.Range("L118").Value = 3919.9998
$ws.Range("M118").Value = -7265
$ws.Range("N118").Value = -6405.9998
$ws.Range("H132").Value = 850.38464
$ws.Range("I132").Value = 607.1429000000001
$ws.Range("K132").Value = 5464.2861
$ws.Range("M132").Value = -2934.2861

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 84273.5
$ws.Range("I113").Value = 143641.72
$ws.Range("J113").Value = 1158
$ws.Range("K113").Value = 143641.72
$ws.Range("L113").Value = 1158
$ws.Range("M113").Value = -141471.72
$ws.Range("N113").Value = -5498
$ws.Range("H138").Value = 58679.4
$ws.Range("J138").Value = 58679.4
$ws.Range("L138").Value = 58679.4
$ws.Range("N138").Value = -68959.39999999999

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 20000500
$ws.Range("J2").Value = 20000500
$ws.Range("L2").Value = 20000500
$ws.Range("N2").Value = -20000724
$ws.Range("H7").Value = 6665.2173
$ws.Range("I7").Value = 9450
$ws.Range("J7").Value = 5180
$ws.Range("K7").Value = 9450
$ws.Range("L7").Value = 5180
$ws.Range("M7").Value = -9338
$ws.Range("N7").Value = -5404
$ws.Range("H46").Value = 1040.3
$ws.Range("I46").Value = 600.25
$ws.Range("J46").Value = 1333.6666
$ws.Range("K46").Value = 600.25
$ws.Range("L46").Value = 1333.6666
$ws.Range("M46").Value = -412.25
$ws.Range("N46").Value = -1709.6666
$ws.Range("H61").Value = 1304
$ws.Range("I61").Value = 1361.6666
$ws.Range("J61").Value = 1246.3334
$ws.Range("K61").Value = 1361.6666
$ws.Range("L61").Value = 1246.3334
$ws.Range("M61").Value = -1159.6666
$ws.Range("N61").Value = -1650.3334
$ws.Range("H113").Value = 1304
$ws.Range("I113").Value = 1361.6666
$ws.Range("J113").Value = 1246.3334
$ws.Range("K113").Value = 1361.6666
$ws.Range("L113").Value = 1246.3334
$ws.Range("M113").Value = 808.3334
$ws.Range("N113").Value = -5586.3334
$ws.Range("H126").Value = 6665.2173
$ws.Range("I126").Value = 9450
$ws.Range("J126").Value = 5180
$ws.Range("K126").Value = 28350
$ws.Range("L126").Value = 15540
$ws.Range("M126").Value = -25880
$ws.Range("N126").Value = -20480
$ws.Range("H136").Value = 14711447
$ws.Range("I136").Value = 25002184
$ws.Range("J136").Value = 10393.571
$ws.Range("K136").Value = 75006552
$ws.Range("L136").Value = 31180.713
$ws.Range("M136").Value = -75004002
$ws.Range("N136").Value = -36280.713
$ws.Range("H139").Value = 46281.273
$ws.Range("J139").Value = 46844.4
$ws.Range("L139").Value = 46844.4
$ws.Range("N139").Value = -57124.4

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 2000
$ws.Range("J5").Value = 2000
$ws.Range("L5").Value = 2000
$ws.Range("N5").Value = -2224
$ws.Range("H138").Value = 59899
$ws.Range("J138").Value = 59899
$ws.Range("L138").Value = 59899
$ws.Range("N138").Value = -70179
